# Auto-generated Excel COM-interop script applying the Kujata_Profits market-price refresh.
$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 372.4375
$ws.Range("I4").Value = 282.7143
$ws.Range("K4").Value = 282.7143
$ws.Range("M4").Value = -168.7143
# Row 17
$ws.Range("H17").Value = 2853
$ws.Range("I17").Value = 453.33334
$ws.Range("J17").Value = 3276.4707
$ws.Range("K17").Value = 1360.00002
$ws.Range("L17").Value = 9829.4121
$ws.Range("M17").Value = -1192.00002
$ws.Range("N17").Value = -10165.4121
# Row 74
$ws.Range("H74").Value = 2900.4
$ws.Range("I74").Value = 2834.3333
$ws.Range("K74").Value = 2834.3333
$ws.Range("M74").Value = -1898.3333
# Row 77
$ws.Range("H77").Value = 2900.4
$ws.Range("I77").Value = 2834.3333
$ws.Range("K77").Value = 14171.6665
$ws.Range("M77").Value = -9491.666499999999
# Row 86
$ws.Range("H86").Value = 2633.4614
$ws.Range("I86").Value = 2791.875
$ws.Range("J86").Value = 2380
$ws.Range("K86").Value = 2791.875
$ws.Range("L86").Value = 2380
$ws.Range("M86").Value = -1668.875
$ws.Range("N86").Value = -4626
# Row 89
$ws.Range("H89").Value = 2633.4614
$ws.Range("I89").Value = 2791.875
$ws.Range("J89").Value = 2380
$ws.Range("K89").Value = 13959.375
$ws.Range("L89").Value = 11900
$ws.Range("M89").Value = -8343.375
$ws.Range("N89").Value = -23132
# Row 118
$ws.Range("H118").Value = 1223.9
$ws.Range("I118").Value = 647.8
$ws.Range("J118").Value = 1800
$ws.Range("K118").Value = 1943.4
$ws.Range("L118").Value = 5400
$ws.Range("M118").Value = -286.3999999999999
$ws.Range("N118").Value = -8714
# Row 137
$ws.Range("H137").Value = 1507.6364
$ws.Range("I137").Value = 1400.5
$ws.Range("J137").Value = 1656.6957
$ws.Range("K137").Value = 4201.5
$ws.Range("L137").Value = 4970.0871
$ws.Range("M137").Value = -1651.5
$ws.Range("N137").Value = -10070.0871
# Row 138
$ws.Range("H138").Value = 1856.09
$ws.Range("I138").Value = 1233.48
$ws.Range("J138").Value = 2063.6267
$ws.Range("K138").Value = 3700.44
$ws.Range("L138").Value = 6190.880099999999
$ws.Range("M138").Value = 1439.56
$ws.Range("N138").Value = -16470.8801

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 899.23334
$ws.Range("I2").Value = 715.6818
$ws.Range("K2").Value = 715.6818
$ws.Range("M2").Value = -602.6818
# Row 32
$ws.Range("H32").Value = 4609.13
$ws.Range("I32").Value = 3795.011
$ws.Range("J32").Value = 13971.5
$ws.Range("K32").Value = 3795.011
$ws.Range("L32").Value = 13971.5
$ws.Range("M32").Value = -3508.011
$ws.Range("N32").Value = -14545.5
# Row 61
$ws.Range("H61").Value = 71429920
$ws.Range("I61").Value = 83334410
$ws.Range("K61").Value = 83334410
$ws.Range("M61").Value = -83334198
# Row 74
$ws.Range("H74").Value = 4179.3
$ws.Range("I74").Value = 3993.3333
$ws.Range("J74").Value = 4259
$ws.Range("K74").Value = 3993.3333
$ws.Range("L74").Value = 4259
$ws.Range("M74").Value = -3119.3333
$ws.Range("N74").Value = -6007
# Row 77
$ws.Range("H77").Value = 4179.3
$ws.Range("I77").Value = 3993.3333
$ws.Range("J77").Value = 4259
$ws.Range("K77").Value = 19966.6665
$ws.Range("L77").Value = 21295
$ws.Range("M77").Value = -15598.6665
$ws.Range("N77").Value = -30031
# Row 116
$ws.Range("H116").Value = 899.23334
$ws.Range("I116").Value = 715.6818
$ws.Range("K116").Value = 715.6818
$ws.Range("M116").Value = 1578.3182
# Row 122
$ws.Range("H122").Value = 2031.4783
$ws.Range("I122").Value = 1889.2778
$ws.Range("J122").Value = 2543.4
$ws.Range("K122").Value = 5667.8334
$ws.Range("L122").Value = 7630.200000000001
$ws.Range("M122").Value = -3217.8334
$ws.Range("N122").Value = -12530.2
# Row 132
$ws.Range("H132").Value = 3436.2258
$ws.Range("I132").Value = 2834.9375
$ws.Range("J132").Value = 4077.6
$ws.Range("K132").Value = 8504.8125
$ws.Range("L132").Value = 12232.8
$ws.Range("M132").Value = -5974.8125
$ws.Range("N132").Value = -17292.8
# Row 136
$ws.Range("H136").Value = 71429920
$ws.Range("I136").Value = 83334410
$ws.Range("K136").Value = 250003230
$ws.Range("M136").Value = -250000680

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 899.23334
$ws.Range("I3").Value = 715.6818
$ws.Range("K3").Value = 715.6818
$ws.Range("M3").Value = -601.6818
# Row 94
$ws.Range("H94").Value = 7812851.5
$ws.Range("I94").Value = 8333632
$ws.Range("J94").Value = 1140
$ws.Range("K94").Value = 8333632
$ws.Range("L94").Value = 1140
$ws.Range("M94").Value = -8333181
$ws.Range("N94").Value = -2042
# Row 134
$ws.Range("H134").Value = 4733.5454
$ws.Range("I134").Value = 7028.3335
$ws.Range("K134").Value = 21085.0005
$ws.Range("M134").Value = -18550.0005

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 2081.5715
$ws.Range("I5").Value = 180.25
$ws.Range("J5").Value = 4616.6665
$ws.Range("K5").Value = 180.25
$ws.Range("L5").Value = 4616.6665
$ws.Range("M5").Value = -68.25
$ws.Range("N5").Value = -4840.6665
# Row 31
$ws.Range("H31").Value = 1275.5238
$ws.Range("I31").Value = 1243.52
$ws.Range("K31").Value = 1243.52
$ws.Range("M31").Value = -948.52
# Row 34
$ws.Range("H34").Value = 1275.5238
$ws.Range("I34").Value = 1243.52
$ws.Range("K34").Value = 1243.52
$ws.Range("M34").Value = -1041.52
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()
# Row 47
$ws.Range("H47").Value = 12000
$ws.Range("J47").Value = 12000
$ws.Range("L47").Value = 12000
$ws.Range("N47").Value = -13132
# Row 100
$ws.Range("H100").Value = 81250
$ws.Range("J100").Value = 81250
$ws.Range("L100").Value = 81250
$ws.Range("N100").Value = -83414
# Row 122
$ws.Range("H122").Value = 790.875
$ws.Range("I122").Value = 805.0769
$ws.Range("J122").Value = 729.3333
$ws.Range("K122").Value = 2415.2307
$ws.Range("L122").Value = 2187.9999
$ws.Range("M122").Value = 34.76929999999993
$ws.Range("N122").Value = -7087.9999

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 636.4318
$ws.Range("I113").Value = 556.7143
$ws.Range("J113").Value = 673.63336
$ws.Range("K113").Value = 1670.1429
$ws.Range("L113").Value = 2020.90008
$ws.Range("M113").Value = 499.8571000000002
$ws.Range("N113").Value = -6360.90008
# Row 136
$ws.Range("H136").Value = 2759.2144
$ws.Range("I136").Value = 2403.75
$ws.Range("J136").Value = 3233.1667
$ws.Range("K136").Value = 7211.25
$ws.Range("L136").Value = 9699.500100000001
$ws.Range("M136").Value = -2111.25
$ws.Range("N136").Value = -19899.5001

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1449.862
$ws.Range("I102").Value = 1495.84
$ws.Range("K102").Value = 1495.84
$ws.Range("M102").Value = 126.1600000000001
# Row 122
$ws.Range("H122").Value = 4144.8184
$ws.Range("I122").Value = 4144.8184
$ws.Range("K122").Value = 12434.4552
$ws.Range("M122").Value = -9984.4552
# Row 126
$ws.Range("H126").Value = 2312.7273
$ws.Range("I126").Value = 1865
$ws.Range("J126").Value = 2568.5715
$ws.Range("K126").Value = 5595
$ws.Range("L126").Value = 7705.7145
$ws.Range("M126").Value = -3125
$ws.Range("N126").Value = -12645.7145
# Row 128
$ws.Range("H128").Value = 37390
# Row 132
$ws.Range("H132").Value = 4661.533
$ws.Range("I132").Value = 5558.7144
$ws.Range("J132").Value = 3876.5
$ws.Range("K132").Value = 16676.1432
$ws.Range("L132").Value = 11629.5
$ws.Range("M132").Value = -14146.1432
$ws.Range("N132").Value = -16689.5
# Row 135
$ws.Range("H135").Value = 35216.43
$ws.Range("J135").Value = 34668.89
$ws.Range("L135").Value = 34668.89
$ws.Range("N135").Value = -44808.89

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1491.5
$ws.Range("I16").Value = 1398.4166
$ws.Range("J16").Value = 2050
$ws.Range("K16").Value = 1398.4166
$ws.Range("L16").Value = 2050
$ws.Range("M16").Value = -1228.4166
$ws.Range("N16").Value = -2390
# Row 40
$ws.Range("H40").Value = 9467.799999999999
$ws.Range("I40").Value = 3744.5
$ws.Range("K40").Value = 3744.5
$ws.Range("M40").Value = -3608.5
# Row 68
$ws.Range("H68").Value = 1243.5714
$ws.Range("I68").Value = 1275.8334
$ws.Range("J68").Value = 1050
$ws.Range("K68").Value = 1275.8334
$ws.Range("L68").Value = 1050
$ws.Range("M68").Value = -526.8334
$ws.Range("N68").Value = -2548
# Row 71
$ws.Range("H71").Value = 1243.5714
$ws.Range("I71").Value = 1275.8334
$ws.Range("J71").Value = 1050
$ws.Range("K71").Value = 6379.166999999999
$ws.Range("L71").Value = 5250
$ws.Range("M71").Value = -2635.166999999999
$ws.Range("N71").Value = -12738
# Row 122
$ws.Range("H122").Value = 22728982
$ws.Range("I122").Value = 35715828
$ws.Range("K122").Value = 107147484
$ws.Range("M122").Value = -107145034

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Range("H103").Value = 12000
$ws.Range("J103").Value = 12000
$ws.Range("L103").Value = 12000
$ws.Range("N103").Value = -14344
# Row 122
$ws.Range("H122").Value = 13890335
$ws.Range("I122").Value = 19232294
$ws.Range("J122").Value = 1241
$ws.Range("K122").Value = 57696882
$ws.Range("L122").Value = 3723
$ws.Range("M122").Value = -57694432
$ws.Range("N122").Value = -8623
# Row 132
$ws.Range("H132").Value = 1809.4286
$ws.Range("I132").Value = 1643.1714
$ws.Range("K132").Value = 4929.5142
$ws.Range("M132").Value = -2399.5142
